$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. MainInfo sheet: change the "SELECT TASKFLOW Nr:" value from 2 to 1.
#    (Dependent formulas in C5/C6 recalc automatically: 10 / 18)
# ---------------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("MainInfo")
$wsMain.Range("C2").Value = 1

# ---------------------------------------------------------------------------
# 2. Results_U3_sub1 (sheet5): update TASKFLOW result text for submissions.
#    Must be written before Results_U1_sub2 so new shared-string entries are
#    created in the same order as in the target workbook.
# ---------------------------------------------------------------------------
$wsR3 = $wb.Worksheets.Item("Results_U3_sub1")
$wsR3.Range("G11").Value = "RESULT(1) TCASE(1)  FLOW(true) MSG(TESTCASE(1):EQUAL`n)"
$wsR3.Range("H11").Value = "RESULT(1) TCASE(2)  FLOW(true) MSG(TESTCASE(2):EQUAL`n)"
$wsR3.Range("G12").Value = "RESULT(2) TCASE(1)  FLOW(true) MSG(TESTCASE(1):EQUAL`n)"
$wsR3.Range("H12").Value = "RESULT(2) TCASE(2)  FLOW(true) MSG(TESTCASE(2):EQUAL`n)"

# ---------------------------------------------------------------------------
# 3. Results_U1_sub2 (sheet3): update TASKFLOW result / error text.
# ---------------------------------------------------------------------------
$wsR1 = $wb.Worksheets.Item("Results_U1_sub2")

$wsR1.Range("G11").Value = "RESULT(1) TCASE(1)  FLOW(true) MSG(TESTCASE(1):EQUAL`n)"
$wsR1.Range("H11").Value = "RESULT(1) TCASE(2)  FLOW(true) MSG(TESTCASE(2):EQUAL`n)"
$wsR1.Range("I11").Value = "RESULT(1) TCASE(3)  FLOW(true) MSG(TESTCASE(3):EQUAL`n)"

$wsR1.Range("G12").Value = "RESULT(2) TCASE(1)  FLOW(true) MSG(TESTCASE(1):NOT-EQUAL`n)"
$wsR1.Range("H12").Value = "RESULT(2) TCASE(2)  FLOW(true) MSG(TESTCASE(2):EQUAL`n)"
$wsR1.Range("I12").Value = "RESULT(2) TCASE(3)  FLOW(true) MSG(TESTCASE(3):EQUAL`n)"
$wsR1.Range("Q12").Value = "ERROR: SUBMIT(2) TESTCASE(1) MSG:(TEXT COMPARE:DEL:(`nCalifornia Poppy+))"

$wsR1.Range("G13").Value = "RESULT(3) TCASE(1)  FLOW(true) MSG(TESTCASE(1):NOT-EQUAL`n)"
$wsR1.Range("H13").Value = "RESULT(3) TCASE(2)  FLOW(true) MSG(TESTCASE(2):EQUAL`n)"
$wsR1.Range("I13").Value = "RESULT(3) TCASE(3)  FLOW(true) MSG(TESTCASE(3):EQUAL`n)"
$wsR1.Range("Q13").Value = "ERROR: SUBMIT(3) TESTCASE(1) MSG:(TEXT COMPARE:DEL:(California Poppy))"

$wsR1.Range("G14").Value = "RESULT(4) TCASE(1)  FLOW(false) MSG(TESTCASE(1):NOT-COMPARED`n)"
$wsR1.Range("H14").Value = "RESULT(4) TCASE(2)  FLOW(false) MSG(TESTCASE(2):NOT-COMPARED`n)"
$wsR1.Range("I14").Value = "RESULT(4) TCASE(3)  FLOW(true) MSG(TESTCASE(3):EQUAL`n)"

# ---------------------------------------------------------------------------
# 4. Fix up sheet selection / active tab: MainInfo becomes the active sheet,
#    Results_U1_sub2 loses its "topLeftCell"/selection and gets a plain
#    single-cell selection at J11 instead.
# ---------------------------------------------------------------------------
$wsR1.Range("J11").Select()
$wsMain.Activate()
